$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 51,5
$data[0,0] = 39400
$data[0,1] = 2007
$data[0,2] = 2.070003986395053
$data[0,3] = 2008
$data[0,4] = -0.3549868696899106
$data[1,0] = 39583
$data[1,1] = 2008
$data[1,2] = 0.5453776865001148
$data[1,3] = 2009
$data[1,4] = -0.7976031984000098
$data[2,0] = 39765
$data[2,1] = 2008
$data[2,2] = 0.517569958955022
$data[2,3] = 2009
$data[2,4] = -5.168396053267498
$data[3,0] = 39948
$data[3,1] = 2009
$data[3,2] = -6.170514117037273
$data[3,3] = 2010
$data[3,4] = -8.396348489509153
$data[4,0] = 40130
$data[4,1] = 2009
$data[4,2] = -3.956152295564896
$data[4,3] = 2010
$data[4,4] = -1.314964327391877
$data[5,0] = 40310
$data[5,1] = 2010
$data[5,2] = -0.2290082001396909
$data[5,3] = 2011
$data[5,4] = -4.327930935900004
$data[6,0] = 40494
$data[6,1] = 2010
$data[6,2] = 1.234995474941392
$data[6,3] = 2011
$data[6,4] = 1.001424185348321
$data[7,0] = 40676
$data[7,1] = 2011
$data[7,2] = 1.406827509327035
$data[7,3] = 2012
$data[7,4] = 2.015050062499957
$data[8,0] = 40862
$data[8,1] = 2011
$data[8,2] = 0.899360810820804
$data[8,3] = 2012
$data[8,4] = 0.475544341751033
$data[9,0] = 41044
$data[9,1] = 2012
$data[9,2] = 1.153683074671208
$data[9,3] = 2013
$data[9,4] = 3.648892256099945
$data[10,0] = 41228
$data[10,1] = 2012
$data[10,2] = 0.9010266119894084
$data[10,3] = 2013
$data[10,4] = 1.506358095275817
$data[11,0] = 41409
$data[11,1] = 2013
$data[11,2] = 0.2186142574756467
$data[11,3] = 2014
$data[11,4] = 0.4006004000999708
$data[12,0] = 41592
$data[12,1] = 2013
$data[12,2] = 0.02019328874804938
$data[12,3] = 2014
$data[12,4] = -1.194807813319188
$data[13,0] = 41774
$data[13,1] = 2014
$data[13,2] = -0.8522658067264599
$data[13,3] = 2015
$data[13,4] = -3.551690943899999
$data[14,0] = 41957
$data[14,1] = 2014
$data[14,2] = 0.1729981757035093
$data[14,3] = 2015
$data[14,4] = 0.6265079396372775
$data[15,0] = 42137
$data[15,1] = 2015
$data[15,2] = -0.2262139320475365
$data[15,3] = 2016
$data[15,4] = -0.7976031983999876
$data[16,0] = 42321
$data[16,1] = 2015
$data[16,2] = 0.09752710595589686
$data[16,3] = 2016
$data[16,4] = -1.022506370243093
$data[17,0] = 42503
$data[17,1] = 2016
$data[17,2] = -0.6258176826215101
$data[17,3] = 2017
$data[17,4] = -0.3994003999000073
$data[18,0] = 42689
$data[18,1] = 2016
$data[18,2] = -0.5280591151586633
$data[18,3] = 2017
$data[18,4] = -0.7240982069264934
$data[19,0] = 42867
$data[19,1] = 2017
$data[19,2] = 0.3239252862367037
$data[19,3] = 2018
$data[19,4] = 1.609625625600009
$data[20,0] = 43053
$data[20,1] = 2017
$data[20,2] = 0.07201851318385799
$data[20,3] = 2018
$data[20,4] = 1.255028673974046
$data[21,0] = 43145
$data[21,1] = 2018
$data[21,2] = 1.456954732048321
$data[21,3] = 2019
$data[21,4] = 2.015050062499957
$data[22,0] = 43235
$data[22,1] = 2018
$data[22,2] = 0.5738128002843901
$data[22,3] = 2019
$data[22,4] = -0.3994003999000184
$data[23,0] = 43326
$data[23,1] = 2018
$data[23,2] = 0.3477859729380528
$data[23,3] = 2019
$data[23,4] = -1.516043567048941
$data[24,0] = 43418
$data[24,1] = 2018
$data[24,2] = 0.3727661260635617
$data[24,3] = 2019
$data[24,4] = -3.305525567352929
$data[25,0] = 43510
$data[25,1] = 2019
$data[25,2] = -0.7761690566734369
$data[25,3] = 2020
$data[25,4] = 0
$data[26,0] = 43600
$data[26,1] = 2019
$data[26,2] = -0.4781004700720293
$data[26,3] = 2020
$data[26,4] = 0.8024032015999882
$data[27,0] = 43691
$data[27,1] = 2019
$data[27,2] = -0.9254001004749823
$data[27,3] = 2020
$data[27,4] = -1.738778148048659
$data[28,0] = 43783
$data[28,1] = 2019
$data[28,2] = -0.801759526476209
$data[28,3] = 2020
$data[28,4] = 1.431264289671219
$data[29,0] = 43875
$data[29,1] = 2020
$data[29,2] = -1.172985875230903
$data[29,3] = 2021
$data[29,4] = -2.378486270400004
$data[30,0] = 43966
$data[30,1] = 2020
$data[30,2] = -1.197849743493773
$data[30,3] = 2021
$data[30,4] = -3.161804390399992
$data[31,0] = 44068
$data[31,1] = 2020
$data[31,2] = -1.503583188367719
$data[31,3] = 2021
$data[31,4] = 5.0514716327553
$data[32,0] = 44159
$data[32,1] = 2020
$data[32,2] = -1.103489789942047
$data[32,3] = 2021
$data[32,4] = 2.294626310579817
$data[33,0] = 44251
$data[33,1] = 2021
$data[33,2] = 3.668278063260222
$data[33,3] = 2022
$data[33,4] = 4.875032525328971
$data[34,0] = 44341
$data[34,1] = 2021
$data[34,2] = 1.064698711638945
$data[34,3] = 2022
$data[34,4] = -2.540956581357878
$data[35,0] = 44432
$data[35,1] = 2021
$data[35,2] = 1.067534122491809
$data[35,3] = 2022
$data[35,4] = 1.042084871410087
$data[36,0] = 44525
$data[36,1] = 2021
$data[36,2] = 0.9704846793491928
$data[36,3] = 2022
$data[36,4] = -0.902682013141165
$data[37,0] = 44617
$data[37,1] = 2022
$data[37,2] = -0.7181024432008964
$data[37,3] = 2023
$data[37,4] = -1.571815848026048
$data[38,0] = 44706
$data[38,1] = 2022
$data[38,2] = -1.524103236349472
$data[38,3] = 2023
$data[38,4] = -1.240907591477092
$data[39,0] = 44798
$data[39,1] = 2022
$data[39,2] = -0.9795431199870586
$data[39,3] = 2023
$data[39,4] = -0.5076503601560978
$data[40,0] = 44890
$data[40,1] = 2022
$data[40,2] = -0.7009264669202708
$data[40,3] = 2023
$data[40,4] = 1.220523709718857
$data[41,0] = 44981
$data[41,1] = 2023
$data[41,2] = 0.331635972555544
$data[41,3] = 2024
$data[41,4] = 1.49657333418427
$data[42,0] = 45071
$data[42,1] = 2023
$data[42,2] = 0.2854413827033664
$data[42,3] = 2024
$data[42,4] = -0.2470349027347551
$data[43,0] = 45163
$data[43,1] = 2023
$data[43,2] = 0.1829021030556488
$data[43,3] = 2024
$data[43,4] = -0.5273501419610804
$data[44,0] = 45254
$data[44,1] = 2023
$data[44,2] = 0.3928252664241905
$data[44,3] = 2024
$data[44,4] = 0.4517021897791018
$data[45,0] = 45345
$data[45,1] = 2024
$data[45,2] = 0.1776394553850924
$data[45,3] = 2025
$data[45,4] = -0.3176932480832284
$data[46,0] = 45436
$data[46,1] = 2024
$data[46,2] = 1.118108578853261
$data[46,3] = 2025
$data[46,4] = 1.532721825047534
$data[47,0] = 45534
$data[47,1] = 2024
$data[47,2] = 0.6979546684258597
$data[47,3] = 2025
$data[47,4] = -0.01999876157223746
$data[48,0] = 45618
$data[48,1] = 2024
$data[48,2] = 0.3224026462283813
$data[48,3] = 2025
$data[48,4] = -0.7185940249203049
$data[49,0] = 45713
$data[49,1] = 2025
$data[49,2] = -2.566037671324872
$data[49,3] = 2026
$data[49,4] = -1.090083898854388
$data[50,0] = 45800
$data[50,1] = 2025
$data[50,2] = -1.490505436658163
$data[50,3] = 2026
$data[50,4] = -0.3349088112516219

$ws.Range("A2:E52").Value = $data

# Add new row 53 (shifts dimension, copy date style format from row above for column A)
$ws.Cells.Item(53,1).Value = 45891
$ws.Cells.Item(53,2).Value = 2025
$ws.Cells.Item(53,3).Value = -2.11737366557071
$ws.Cells.Item(53,4).Value = 2026
$ws.Cells.Item(53,5).Value = -0.5919451648311758
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Cells.Item(53,1).Value = 45891
